$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-8 from 45185 (2023-09-16)
# to 45204 (2023-10-05), keeping the existing date number format.
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45204
}
